# DDT using Excel File & Cucumbar simple report
# Populate the "EmployeeLoginCredentials" sheet with the same login-credentials
# table already present on the "Login" sheet, and make it the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EmployeeLoginCredentials")
$ws2 = $wb.Worksheets.Item("Login")

# ---- Copy the data table from Login into EmployeeLoginCredentials ----
$headers = @("FirstName", "LastName", "UserName", "password")
for ($c = 1; $c -le 4; $c++) {
    $ws1.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$data = @(
    @("John", "Smith", "John5623", "AmirKhan_@123"),
    @("Mary", "Ann",   "Mary5612", "AmirKhan_@124"),
    @("Ali",  "Din",   "Ali5623",  "AmirKhan_@124")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---- Match the "Login" sheet's look (font sizes) ----
$ws1.Range("A1:D4").Font.Size = 20
$ws1.Range("D2:D4").Font.Name = "Menlo"
$ws1.Range("D2:D4").Font.Color = 16750080

# ---- Column widths (characters), matching the Login sheet's proportions ----
$ws1.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws1.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws1.Columns.Item(3).ColumnWidth = 20.666666666666668
$ws1.Columns.Item(4).ColumnWidth = 33.666666666666664

# ---- Make EmployeeLoginCredentials the active sheet/tab ----
$ws1.Activate()
$ws1.Range("A1:C1048576").Select()

# ---- Login sheet selection now spans the full sheet ----
$ws2.Range("A1:XFD1048576").Select()
